$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDR_B")

# --- Row 4 / Row 5 : DQS1 pair re-measured (names swap, values updated) ---
$ws.Range("B4").Value = "DRAM_DQS1_B_P"
$ws.Range("E4").Value = 18.53
$ws.Range("G4").Value = 20.93
$ws.Range("H4").Value = 108.4
$ws.Range("J4").Value = 55.5
$ws.Range("L4").Value = 199.9

$ws.Range("B5").Value = "DRAM_DQS1_B_N"
$ws.Range("E5").Value = 18.55
$ws.Range("G5").Value = 20.95
$ws.Range("H5").Value = 108.67
$ws.Range("J5").Value = 55
$ws.Range("L5").Value = 199.67

# --- Row 14 : DRAM_D15_B re-measured ---
$ws.Range("E14").Value = 21.84
$ws.Range("G14").Value = 22.24
$ws.Range("H14").Value = 150.9
$ws.Range("L14").Value = 199.7

# --- Row 15 / Row 16 : DQS0 pair re-measured (names swap, values updated) ---
$ws.Range("B15").Value = "DRAM_DQS0_B_P"
$ws.Range("E15").Value = 10.73
$ws.Range("G15").Value = 13.13
$ws.Range("H15").Value = 62.83
$ws.Range("L15").Value = 147.43

$ws.Range("B16").Value = "DRAM_DQS0_B_N"
$ws.Range("E16").Value = 10.71
$ws.Range("G16").Value = 13.11
$ws.Range("H16").Value = 62.74
$ws.Range("L16").Value = 147.34

# --- Rows 17-25 : byte lane 0 nets reordered/re-measured to fix DRC errors ---
$ws.Range("B17").Value = "DRAM_D05_B"
$ws.Range("E17").Value = 16.7
$ws.Range("G17").Value = 16.7
$ws.Range("H17").Value = 97.62
$ws.Range("J17").Value = 49.8
$ws.Range("L17").Value = 147.42

$ws.Range("B18").Value = "DRAM_DMI0_B"
$ws.Range("E18").Value = 16.17
$ws.Range("G18").Value = 16.17
$ws.Range("H18").Value = 94.58
$ws.Range("J18").Value = 52.8
$ws.Range("L18").Value = 147.38

$ws.Range("B19").Value = "DRAM_D04_B"
$ws.Range("E19").Value = 16.05
$ws.Range("G19").Value = 16.05
$ws.Range("H19").Value = 93.74
$ws.Range("J19").Value = 53.6
$ws.Range("L19").Value = 147.34

$ws.Range("B20").Value = "DRAM_D03_B"
$ws.Range("E20").Value = 18.04
$ws.Range("G20").Value = 18.04
$ws.Range("H20").Value = 105.34
$ws.Range("J20").Value = 42

$ws.Range("B21").Value = "DRAM_D06_B"
$ws.Range("E21").Value = 15.84
$ws.Range("G21").Value = 15.84
$ws.Range("H21").Value = 92.61
$ws.Range("J21").Value = 54.7
$ws.Range("L21").Value = 147.31

$ws.Range("B22").Value = "DRAM_D07_B"
$ws.Range("E22").Value = 16.99
$ws.Range("G22").Value = 16.99
$ws.Range("H22").Value = 99.27
$ws.Range("J22").Value = 48
$ws.Range("L22").Value = 147.27

$ws.Range("B23").Value = "DRAM_D00_B"
$ws.Range("E23").Value = 16.41
$ws.Range("G23").Value = 16.41
$ws.Range("H23").Value = 95.87
$ws.Range("J23").Value = 51.5
$ws.Range("L23").Value = 147.37

$ws.Range("B24").Value = "DRAM_D02_B"
$ws.Range("E24").Value = 15.84
$ws.Range("G24").Value = 15.84
$ws.Range("H24").Value = 92.63
$ws.Range("J24").Value = 54.5

$ws.Range("B25").Value = "DRAM_D01_B"
$ws.Range("E25").Value = 16.64
$ws.Range("G25").Value = 16.64
$ws.Range("H25").Value = 97.23
$ws.Range("J25").Value = 49.9
$ws.Range("L25").Value = 147.13

# --- Leave behind the orphaned conditional-format dxf (white fill) that
#     Excel keeps in styles.xml once a cell-highlight rule used while
#     checking the fix is cleared again ---
$cfRange = $ws.Range("B4")
$cf = $cfRange.FormatConditions.Add(1, 3, "1")
$cf.Interior.Color = 16777215
$cfRange.FormatConditions.Delete()

# --- Update the saved selection / active cell for the sheet ---
$ws.Activate()
$ws.Range("F38").Select() | Out-Null
